# Update gh-pages output data (想去人数 / "want to go" counters) for the
# 苏州-漫展信息 workbook, as published at commit 456a3b4.
#
# Sheet "展览" (index 1) and sheet "全部类型" (index 4) both list the same
# events (the latter aggregates rows from every other sheet), so the same
# logical updates are applied to column F ("想去人数") on both sheets -
# just at different row numbers.

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item(1)   # 展览
$sheetAllTypes    = $wb.Worksheets.Item(4)  # 全部类型

# Row -> new value updates for the "展览" sheet
$exhibitionUpdates = @{
    "F2"  = 1036
    "F3"  = 13631
    "F4"  = 41
    "F7"  = 1746
    "F8"  = 155
    "F9"  = 125
    "F10" = 86
    "F12" = 4
    "F14" = 13624
    "F16" = 606
    "F17" = 8994
    "F18" = 10
    "F19" = 8098
    "F20" = 257
    "F21" = 13
    "F30" = 396
    "F32" = 195
    "F33" = 381
    "F35" = 10
}

foreach ($cellRef in $exhibitionUpdates.Keys) {
    $sheetExhibition.Range($cellRef).Value = $exhibitionUpdates[$cellRef]
}

# Row -> new value updates for the "全部类型" sheet
$allTypesUpdates = @{
    "F2"  = 1036
    "F3"  = 13631
    "F4"  = 41
    "F7"  = 1746
    "F8"  = 155
    "F9"  = 125
    "F10" = 86
    "F12" = 4
    "F14" = 13624
    "F16" = 606
    "F17" = 8994
    "F18" = 10
    "F19" = 8098
    "F20" = 257
    "F21" = 13
    "F32" = 396
    "F34" = 195
    "F35" = 381
    "F37" = 10
}

foreach ($cellRef in $allTypesUpdates.Keys) {
    $sheetAllTypes.Range($cellRef).Value = $allTypesUpdates[$cellRef]
}
